$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update text labels (shorten descriptions, dropping "Chapa ..." prefixes)
$ws.Range("B3").Value = "Negra"
$ws.Range("B4").Value = "Galvanizada"
$ws.Range("B6").Value = "Semilla de melón"

# Adjust row heights (Mac -> Windows font metrics tweak)
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 18.75
$ws.Rows.Item(6).RowHeight = 18.75

$wb.Save()
